$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (F and G) before the existing "Is Significant" column,
# which pushes "Is Significant" (and its data) from column F to column H.
# Excel COM automatically extends the formatting (style) of the surrounding
# header cells to the newly inserted cells.
$ws.Range("F1:G5").EntireColumn.Insert()

# New header cells
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# New "Observed" values (column F)
$ws.Range("F2").Value = "[393 577] ; [16  4]"
$ws.Range("F3").Value = "[211 374] ; [15  5]"
$ws.Range("F4").Value = "[211 374] ; [15  5]"
$ws.Range("F5").Value = "[455 518] ; [17  3]"

# New "Expected" values (column G)
$ws.Range("G2").Value = "[400.73737374 569.26262626] ; [ 8.26262626 11.73737374]"
$ws.Range("G3").Value = "[218.52892562 366.47107438] ; [ 7.47107438 12.52892562]"
$ws.Range("G4").Value = "[218.52892562 366.47107438] ; [ 7.47107438 12.52892562]"
$ws.Range("G5").Value = "[462.49345418 510.50654582] ; [ 9.50654582 10.49345418]"
